{"js": "// Oct 30 Lab Updated \u2014 \"Now includes an output function\"\n//\n// The canonical diff shows the \"_GoBack\" bookmark (Word's \"last edit\n// position\" marker) moving from the very end of the document to a spot\n// in the middle of the word \"header\" inside the sentence:\n//   \"#include \"rsfunct.h\" rsfunct is my custom header that prototypes\n//    functions for the average daily balance calculator\"\n// i.e. the run is split into \"...custom head\" + [bookmark] + \"er that\n// prototypes...\". The visible text is unchanged; only the bookmark\n// position (and the run split it implies) changes.\n\n// 1) Remove the pre-existing \"_GoBack\" bookmark (it currently sits at the\n//    very end of the document, right after \"...without user interaction\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Locate the sentence fragment that ends right where the new bookmark\n//    belongs (\"...custom head\") and collapse a zero-width range right\n//    after it.\nconst searchResults = context.document.body.search(\n  \"#include \\u201crsfunct.h\\u201d rsfunct is my custom head\",\n  { matchCase: true }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  const hit = searchResults.items[0];\n  const insertionPoint = hit.getRange(\"After\");\n\n  // 3) Drop the \"_GoBack\" bookmark at that split point (this is what\n  //    forces the surrounding run to split into two runs on save).\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Oct 30 Lab Updated \u2014 \"Now includes an output function\"\n#\n# The canonical diff shows the \"_GoBack\" bookmark (Word's \"last edit\n# position\" marker) moving from the very end of the document to a spot\n# in the middle of the word \"header\" inside the sentence:\n#   #include \"rsfunct.h\" rsfunct is my custom header that prototypes\n#   functions for the average daily balance calculator\n# i.e. the run that held that sentence gets split into \"...custom head\"\n# + [bookmark] + \"er that prototypes...\". The visible text is unchanged;\n# only the bookmark position (and the run split it implies) changes.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the pre-existing \"_GoBack\" bookmark \u2014 it currently sits at the\n#    very end of the document, right after \"...without user interaction\".\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Find the sentence fragment that ends right where the new bookmark\n#    belongs (\"...custom head\") and collapse the found range to its end.\n$quoteOpen = [char]8220\n$quoteClose = [char]8221\n$target = $quoteOpen + \"rsfunct.h\" + $quoteClose + \" rsfunct is my custom head\"\n\n$rng = $d.Content\n$found = $rng.Find.Execute($target)\n\nif ($found) {\n    $rng.Collapse(0)  # wdCollapseEnd\n\n    # 3) Drop the \"_GoBack\" bookmark at that split point (this is what\n    #    forces the surrounding run to split into two runs on save).\n    $d.Bookmarks.Add(\"_GoBack\", $rng)\n}\n"}
